$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 29745
$ws.Range("J13").Value = 29745
$ws.Range("L13").Value = 29745
$ws.Range("N13").Value = -30083

$ws.Range("H54").Value = 40999.668
$ws.Range("I54").Value = 22999.666
$ws.Range("J54").Value = 49999.668
$ws.Range("K54").Value = 22999.666
$ws.Range("L54").Value = 49999.668
$ws.Range("M54").Value = -22513.666
$ws.Range("N54").Value = -50971.668

$ws.Range("H69").Value = 4000
$ws.Range("I69").Value = 4000
$ws.Range("K69").Value = 12000
$ws.Range("M69").Value = -11126

$ws.Range("H70").Value = 1945.8695
$ws.Range("I70").Value = 1484.7
$ws.Range("J70").Value = 2300.6155
$ws.Range("K70").Value = 4454.1
$ws.Range("L70").Value = 6901.8465
$ws.Range("M70").Value = -4184.1
$ws.Range("N70").Value = -7441.8465

$ws.Range("H72").Value = 4000
$ws.Range("I72").Value = 4000
$ws.Range("K72").Value = 36000
$ws.Range("M72").Value = -31632

$ws.Range("H73").Value = 1945.8695
$ws.Range("I73").Value = 1484.7
$ws.Range("J73").Value = 2300.6155
$ws.Range("K73").Value = 4454.1
$ws.Range("L73").Value = 6901.8465
$ws.Range("M73").Value = -3518.1
$ws.Range("N73").Value = -8773.8465

$ws.Range("H76").Value = 4047.4736
$ws.Range("I76").Value = 3193.2
$ws.Range("J76").Value = 7251
$ws.Range("K76").Value = 3193.2
$ws.Range("L76").Value = 7251
$ws.Range("M76").Value = -2878.2
$ws.Range("N76").Value = -7881

$ws.Range("H79").Value = 4047.4736
$ws.Range("I79").Value = 3193.2
$ws.Range("J79").Value = 7251
$ws.Range("K79").Value = 3193.2
$ws.Range("L79").Value = 7251
$ws.Range("M79").Value = -2101.2
$ws.Range("N79").Value = -9435

$ws.Range("H86").Value = 1356.2858
$ws.Range("I86").Value = 1098.8
$ws.Range("K86").Value = 1098.8
$ws.Range("M86").Value = 24.20000000000005

$ws.Range("H89").Value = 1356.2858
$ws.Range("I89").Value = 1098.8
$ws.Range("K89").Value = 5494
$ws.Range("M89").Value = 122

$ws.Range("H92").Value = 1441.2222
$ws.Range("I92").Value = 1138.3572
$ws.Range("K92").Value = 1138.3572
$ws.Range("M92").Value = 109.6428000000001

$ws.Range("H93").Value = 27188.354
$ws.Range("J93").Value = 27188.354
$ws.Range("L93").Value = 27188.354
$ws.Range("N93").Value = -32180.354

$ws.Range("H137").Value = 1702688.9
$ws.Range("I137").Value = 2268675.5
$ws.Range("K137").Value = 6806026.5
$ws.Range("M137").Value = -6803476.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6040.2
$ws.Range("I74").Value = 9061
$ws.Range("J74").Value = 3019.4
$ws.Range("K74").Value = 9061
$ws.Range("L74").Value = 3019.4
$ws.Range("M74").Value = -8187
$ws.Range("N74").Value = -4767.4

$ws.Range("H77").Value = 6040.2
$ws.Range("I77").Value = 9061
$ws.Range("J77").Value = 3019.4
$ws.Range("K77").Value = 45305
$ws.Range("L77").Value = 15097
$ws.Range("M77").Value = -40937
$ws.Range("N77").Value = -23833

$ws.Range("H132").Value = 2232.3462
$ws.Range("I132").Value = 978.7273
$ws.Range("J132").Value = 3151.6667
$ws.Range("K132").Value = 2936.1819
$ws.Range("L132").Value = 9455.000100000001
$ws.Range("M132").Value = -406.1819
$ws.Range("N132").Value = -14515.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 36499.75
$ws.Range("J59").Value = 36499.75
$ws.Range("L59").Value = 36499.75
$ws.Range("N59").Value = -38193.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2775.6287
$ws.Range("I31").Value = 984.3889
$ws.Range("K31").Value = 984.3889
$ws.Range("M31").Value = -689.3889

$ws.Range("H34").Value = 2775.6287
$ws.Range("I34").Value = 984.3889
$ws.Range("K34").Value = 984.3889
$ws.Range("M34").Value = -782.3889

$ws.Range("H132").Value = 2252.5833
$ws.Range("I132").Value = 1288.762
$ws.Range("K132").Value = 3866.286
$ws.Range("M132").Value = -1336.286

$ws.Range("H135").Value = 39926.668
$ws.Range("J135").Value = 39926.668
$ws.Range("L135").Value = 39926.668
$ws.Range("N135").Value = -50066.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4630239
$ws.Range("I113").Value = 641.6923
$ws.Range("J113").Value = 8929151
$ws.Range("K113").Value = 1925.0769
$ws.Range("L113").Value = 26787453
$ws.Range("M113").Value = 244.9231
$ws.Range("N113").Value = -26791793

$ws.Range("H124").Value = 4347.1816
$ws.Range("I124").Value = 2227.5
$ws.Range("J124").Value = 9999.666999999999
$ws.Range("K124").Value = 6682.5
$ws.Range("L124").Value = 29999.001
$ws.Range("M124").Value = -1772.5
$ws.Range("N124").Value = -39819.001

$ws.Range("H134").Value = 4475
$ws.Range("I134").Value = 4129.1665
$ws.Range("J134").Value = 4890
$ws.Range("K134").Value = 12387.4995
$ws.Range("L134").Value = 14670
$ws.Range("M134").Value = -7317.499500000002
$ws.Range("N134").Value = -24810

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2761.0625
$ws.Range("I132").Value = 1602.85
$ws.Range("J132").Value = 4691.4165
$ws.Range("K132").Value = 4808.549999999999
$ws.Range("L132").Value = 14074.2495
$ws.Range("M132").Value = -2278.549999999999
$ws.Range("N132").Value = -19134.2495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 27890
$ws.Range("J121").Value = 27890
$ws.Range("L121").Value = 27890
$ws.Range("N121").Value = -31384

$ws.Range("H132").Value = 5491.9165
$ws.Range("I132").Value = 1168.7
$ws.Range("K132").Value = 3506.1
$ws.Range("M132").Value = -976.1000000000004

$ws.Range("H136").Value = 4235.75
$ws.Range("I136").Value = 1088.1666
$ws.Range("J136").Value = 7383.3335
$ws.Range("K136").Value = 3264.4998
$ws.Range("L136").Value = 22150.0005
$ws.Range("M136").Value = -714.4998000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 30550
$ws.Range("J92").Value = 30550
$ws.Range("L92").Value = 30550
$ws.Range("N92").Value = -35542

$ws.Range("H121").Value = 28890
$ws.Range("J121").Value = 28890
$ws.Range("L121").Value = 28890
$ws.Range("N121").Value = -32384

$ws.Range("H132").Value = 6805211.5
$ws.Range("I132").Value = 1623
$ws.Range("J132").Value = 15155070
$ws.Range("K132").Value = 4869
$ws.Range("L132").Value = 45465210
$ws.Range("M132").Value = -2339

$ws.Range("H136").Value = 4314.12
$ws.Range("J136").Value = 9000.625
$ws.Range("L136").Value = 27001.875
$ws.Range("N136").Value = -32101.875
